$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions data pull).
# D-column price strings are digit-only-looking text (e.g. "1.00", "70.877.95")
# that Excel's COM layer would otherwise auto-coerce to a Number cell, losing
# the exact text formatting (trailing zeros, thousand-dot grouping) and the
# inline-string cell type. Force text via a temporary "@" number format, then
# ClearFormats() so the cell keeps its original (unstyled) appearance.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '70.877.95'
$ws.Range('E2').Value = '  -0.24%  '
Set-TextValue $ws.Range('D3') '3.843.99'
$ws.Range('E3').Value = '  +1.11%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws.Range('D5') '702.76'
$ws.Range('E5').Value = '  -0.40%  '
Set-TextValue $ws.Range('D6') '172.30'
$ws.Range('E6').Value = '  -0.23%  '
Set-TextValue $ws.Range('D7') '3.842.55'
$ws.Range('E7').Value = '  +1.11%  '
Set-TextValue $ws.Range('D8') '1.00'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').Value = '  -1.09%  '
Set-TextValue $ws.Range('D11') '7.30'
$ws.Range('E11').Value = '  -2.51%  '
Set-TextValue $ws.Range('D12') '0.458'
$ws.Range('E12').Value = '  -0.80%  '
Set-TextValue $ws.Range('D13') '0.0000255'
$ws.Range('E13').Value = '  -2.36%  '
Set-TextValue $ws.Range('D14') '36.28'
$ws.Range('E14').Value = '  +0.16%  '
Set-TextValue $ws.Range('D15') '4.492.56'
$ws.Range('E15').Value = '  +1.11%  '
Set-TextValue $ws.Range('D16') '3.902.38'
$ws.Range('E16').Value = '  +2.59%  '
Set-TextValue $ws.Range('D17') '70.970.85'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('E19').Value = '  +0.59%  '
Set-TextValue $ws.Range('D20') '17.38'
$ws.Range('E20').Value = '  -2.90%  '
Set-TextValue $ws.Range('D21') '10.74'
$ws.Range('E21').Value = '  -4.33%  '
Set-TextValue $ws.Range('D22') '491.82'
$ws.Range('E22').Value = '  +1.52%  '
Set-TextValue $ws.Range('D23') '0.716'
$ws.Range('E23').Value = '  +0.03%  '
Set-TextValue $ws.Range('D24') '84.87'
$ws.Range('E24').Value = '  +1.45%  '
Set-TextValue $ws.Range('D25') '0.0000146'
$ws.Range('E25').Value = '  -0.71%  '
Set-TextValue $ws.Range('D26') '12.11'
$ws.Range('E26').Value = '  -2.17%  '
Set-TextValue $ws.Range('D27') '10.49'
$ws.Range('E27').Value = '  -0.69%  '
Set-TextValue $ws.Range('D28') '2.11'
$ws.Range('E28').Value = '  -3.12%  '
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('E30').Value = '  +0.04%  '
Set-TextValue $ws.Range('D31') '7.48'
$ws.Range('E31').Value = '  -1.05%  '
Set-TextValue $ws.Range('D32') '2.26'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('E33').Value = '  +2.15%  '
Set-TextValue $ws.Range('D34') '29.33'
$ws.Range('E34').Value = '  -0.89%  '
Set-TextValue $ws.Range('D35') '3.800.80'
$ws.Range('E35').Value = '  +1.22%  '
Set-TextValue $ws.Range('D36') '9.12'
$ws.Range('E36').Value = '  -1.24%  '
Set-TextValue $ws.Range('D37') '1.00'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -0.22%  '
Set-TextValue $ws.Range('D39') '2.36'
$ws.Range('E39').Value = '  +6.15%  '
$ws.Range('E40').Value = '  +7.04%  '
Set-TextValue $ws.Range('D41') '6.00'
$ws.Range('E41').Value = '  +0.40%  '
Set-TextValue $ws.Range('D42') '3.29'
$ws.Range('E42').Value = '  -6.00%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  +0.11%  '
Set-TextValue $ws.Range('D45') '0.000311'
$ws.Range('E45').Value = '  -5.47%  '
Set-TextValue $ws.Range('D46') '163.72'
$ws.Range('E46').Value = '  +1.25%  '
Set-TextValue $ws.Range('D47') '48.71'
$ws.Range('E47').Value = '  -1.45%  '
Set-TextValue $ws.Range('D48') '0.298'
$ws.Range('E48').Value = '  -0.83%  '
Set-TextValue $ws.Range('D49') '8.62'
$ws.Range('E49').Value = '  +0.67%  '
Set-TextValue $ws.Range('D50') '43.37'
$ws.Range('E50').Value = '  -3.89%  '
Set-TextValue $ws.Range('D51') '408.53'
$ws.Range('E51').Value = '  +2.39%  '
